$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FRED Graph")

# Updated GDP observations (revised data + new algorithm) for rows 12-30, column B
$ws.Cells.Item(12, 2).Value = 262794.28700000001
$ws.Cells.Item(13, 2).Value = 279284.49900000001
$ws.Cells.Item(14, 2).Value = 296612.614
$ws.Cells.Item(15, 2).Value = 322239.14500000002
$ws.Cells.Item(16, 2).Value = 346494.92599999998
$ws.Cells.Item(17, 2).Value = 362080.02799999999
$ws.Cells.Item(18, 2).Value = 378592.255
$ws.Cells.Item(19, 2).Value = 397776.821
$ws.Cells.Item(20, 2).Value = 403211.91100000002
$ws.Cells.Item(21, 2).Value = 425450.39399999997
$ws.Cells.Item(22, 2).Value = 438793.74300000002
$ws.Cells.Item(23, 2).Value = 448273.18400000001
$ws.Cells.Item(24, 2).Value = 455796.80499999999
$ws.Cells.Item(25, 2).Value = 469122.70699999999
$ws.Cells.Item(26, 2).Value = 488987.413
$ws.Cells.Item(27, 2).Value = 508569.36
$ws.Cells.Item(28, 2).Value = 525176.09299999999
$ws.Cells.Item(29, 2).Value = 547123.10900000005
$ws.Cells.Item(30, 2).Value = 566892.03799999994

# Add new row 31 for the 2020-01-01 observation.
# Copy formatting (date format in A, numeric format in B) from row 30 first,
# then set the new values.
$ws.Range("A30:B30").Copy()
$ws.Range("A31:B31").PasteSpecial(-4122)
$ws.Cells.Item(31, 1).Value = 43831
$ws.Cells.Item(31, 2).Value = 561027.94099999999

$ws.Application.CutCopyMode = $false

# Dimension / extent grows to A1:B31 automatically; update the active selection
# to B1 to match the saved workbook view.
$ws.Range("B1").Select()
